$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    if ($v -eq 45171) {
        $cell.Value = 45172
    }
}
